# Apply cryptos list update (coin prices / 1h volume %) per commit.
#
# Note: for columns whose new value looks like a plain number (e.g. "227.87"),
# the literal value below starts with a leading apostrophe ('). That mirrors
# what a user typing into Excel would do to force the cell to stay literal
# text (matching the source workbook's inlineStr cells) instead of being
# silently parsed/reformatted as a numeric value (which would turn
# "104.30" into 104.3, "69.90" into 69.9, etc). The apostrophe itself is
# not stored in the cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''35.491.81'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = '''1.847.19'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''227.87'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '''41.21'
$ws.Range('E8').Value = '  +7.67%  '
$ws.Range('D9').Value = '''0.307'
$ws.Range('E9').Value = '  +5.09%  '
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  +3.47%  '
$ws.Range('D12').Value = '''2.116.63'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').Value = '''11.59'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '''1.846.60'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '''4.73'
$ws.Range('E15').Value = '  +6.36%  '
$ws.Range('D16').Value = '''0.666'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').Value = '''35.416.13'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').Value = '''69.90'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '''245.31'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '0.0₃0793'
$ws.Range('D21').Value = '''12.19'
$ws.Range('E21').Value = '  +8.03%  '
$ws.Range('D22').Value = '''4.78'
$ws.Range('E22').Value = '  +15.42%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '''2.20'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').Value = '''170.91'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = '''7.89'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '''17.84'
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('D28').Value = '''0.123'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('D29').Value = '''3.460.74'
$ws.Range('E29').Value = '  +42.44%  '
$ws.Range('E31').Value = '  +7.81%  '
$ws.Range('D32').Value = '''3.93'
$ws.Range('E32').Value = '  +2.89%  '
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('E34').Value = '  +1.97%  '
$ws.Range('E35').Value = '  +2.81%  '
$ws.Range('D36').Value = '''0.677'
$ws.Range('E37').Value = '  +9.54%  '
$ws.Range('D38').Value = '''88.66'
$ws.Range('E38').Value = '  +9.17%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '''1.338.69'
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.08'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('D41').Value = '''0.0194'
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('D42').Value = '''2.41'
$ws.Range('E42').Value = '  +2.43%  '
$ws.Range('D43').Value = '''1.28'
$ws.Range('E43').Value = '  +4.54%  '
$ws.Range('D44').Value = '''14.82'
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('E47').Value = '  +3.68%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '''6.06'
$ws.Range('E48').Value = '  +4.82%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '''2.015.15'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('D50').Value = '''104.30'
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('E51').Value = '  +0.12%  '
